# Updates the "cryptos" list on Sheet1 with refreshed price / 1h-volume
# data (as produced by the periodic GitHub Actions scraper job), and
# re-orders a couple of rows whose ranking changed (USDe/PEPE swapped
# positions 30/31, and WhiteBITCoin/RenderToken/Bittensor rotated at
# positions 45/46/47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric ("5.47", "0.994", "1.00", ...) must
# be forced to Text format *before* the value is assigned. Otherwise the
# COM layer auto-converts the string into a real number (e.g. "0.400"
# would become 0.4), which would not match the source workbook where
# every cell in this table is stored as text.
$numericTextCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D12", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "54.795.96"
$ws.Range("E2").Value = "  +9.36%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.426.17"
$ws.Range("E3").Value = "  +10.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.45%  "

# Row 5 - BNB
$ws.Range("D5").Value = "476.38"
$ws.Range("E5").Value = "  +14.52%  "

# Row 6 - Solana
$ws.Range("D6").Value = "139.59"
$ws.Range("E6").Value = "  +23.75%  "

# Row 7 - USDC
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.55%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.502"
$ws.Range("E8").Value = "  +14.96%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.452.55"
$ws.Range("E9").Value = "  +12.51%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0957"
$ws.Range("E10").Value = "  +14.99%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  +12.83%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  +13.54%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +4.18%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.836.92"
$ws.Range("E14").Value = "  +9.68%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "54.853.55"
$ws.Range("E15").Value = "  +9.90%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "20.35"
$ws.Range("E16").Value = "  +14.50%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  +20.45%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.435.35"
$ws.Range("E18").Value = "  +9.99%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "4.33"
$ws.Range("E19").Value = "  +13.07%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "9.93"
$ws.Range("E20").Value = "  +20.57%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "312.29"
$ws.Range("E21").Value = "  +10.72%  "

# Row 22 - Dai
$ws.Range("D22").Value = "0.995"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.66"
$ws.Range("E23").Value = "  +16.56%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "57.02"
$ws.Range("E24").Value = "  +12.48%  "

# Row 25 - Binance-PegBSC-USD
$ws.Range("E25").Value = "  +0.36%  "

# Row 26 - Polygon
$ws.Range("D26").Value = "0.400"
$ws.Range("E26").Value = "  +13.95%  "

# Row 27 - Kaspa
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +27.40%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.540.98"
$ws.Range("E28").Value = "  +9.70%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.28"
$ws.Range("E29").Value = "  +13.74%  "

# Row 30 - was USDe, now PEPE (rank swap with row 31)
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0769"
$ws.Range("E30").Value = "  +25.79%  "

# Row 31 - was PEPE, now USDe (rank swap with row 30)
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 - Monero
$ws.Range("D32").Value = "147.95"
$ws.Range("E32").Value = "  +5.99%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "18.08"
$ws.Range("E33").Value = "  +13.70%  "

# Row 34 - PancakeSwap
$ws.Range("D34").Value = "1.46"
$ws.Range("E34").Value = "  +15.78%  "

# Row 35 - Aptos
$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  +14.43%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "1.12"
$ws.Range("E36").Value = "  +20.42%  "

# Row 37 - NEARProtocol
$ws.Range("D37").Value = "3.57"
$ws.Range("E37").Value = "  +14.45%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "0.831"
$ws.Range("E38").Value = "  +16.57%  "

# Row 39 - OKB
$ws.Range("D39").Value = "33.63"
$ws.Range("E39").Value = "  +7.66%  "

# Row 40 - FirstDigitalUSD
$ws.Range("D40").Value = "0.988"
$ws.Range("E40").Value = "  +0.28%  "

# Row 41 - Filecoin
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  +14.96%  "

# Row 42 - Hedera
$ws.Range("D42").Value = "0.0543"
$ws.Range("E42").Value = "  +15.41%  "

# Row 43 - Mantle
$ws.Range("D43").Value = "0.589"
$ws.Range("E43").Value = "  +12.76%  "

# Row 44 - Stacks
$ws.Range("D44").Value = "1.28"
$ws.Range("E44").Value = "  +18.33%  "

# Row 45 - was WhiteBITCoin, now RenderToken (3-way rotation with rows 46/47)
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.68"
$ws.Range("E45").Value = "  +32.08%  "

# Row 46 - was RenderToken, now Bittensor (3-way rotation with rows 45/47)
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "256.42"
$ws.Range("E46").Value = "  +39.20%  "

# Row 47 - was Bittensor, now WhiteBITCoin (3-way rotation with rows 45/46)
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "10.12"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48 - Stellar
$ws.Range("D48").Value = "0.0889"
$ws.Range("E48").Value = "  +16.18%  "

# Row 49 - VeChain
$ws.Range("D49").Value = "0.0221"
$ws.Range("E49").Value = "  +14.72%  "

# Row 50 - Maker
$ws.Range("D50").Value = "1.891.85"
$ws.Range("E50").Value = "  +5.11%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "16.85"
$ws.Range("E51").Value = "  +14.10%  "
